$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text updates (rich-text runs: replace only the changed substring so
#     the surrounding run formatting is left intact) ---

# "Volume 30   Number  7" -> "Volume 30   Number  9"
$volRange = $ws.Range("A8")
$volText = $volRange.Text
$volIdx = $volText.LastIndexOf("7")
$volRange.Characters($volIdx + 1, 1).Text = "9"

# "Report Covering the Week  2/13/2023  Through  2/19/2023"
#   -> "Report Covering the Week  2/27/2023  Through  3/5/2023"
$weekRange = $ws.Range("C9")
$weekText = $weekRange.Text
$idx1 = $weekText.IndexOf("2/13/2023")
$weekRange.Characters($idx1 + 1, 9).Text = "2/27/2023"
$weekText2 = $ws.Range("C9").Text
$idx2 = $weekText2.IndexOf("2/19/2023")
$ws.Range("C9").Characters($idx2 + 1, 9).Text = "3/5/2023"

# --- Cells that change data type (numeric <-> text placeholder) ---
# These reuse an already-formatted donor cell (C14, style index 14) to restore the
# exact "right aligned / General" placeholder look after forcing text, and a
# "#,##0" number format when turning a placeholder back into a real number.

# F14: placeholder "0" -> real number 2
$ws.Range("F14").Value = 2
$ws.Range("F14").NumberFormat = "#,##0"

# C30: placeholder "0" -> real number 1
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"

# D28 / E28: real numbers -> placeholders "0" / "***.*"
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# D29 / E29: real numbers -> placeholders "0" / "***.*"
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

# --- Bulk weekly crime-statistics refresh (rows 14-30, columns C:N) ---
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -33.333333333333
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = -33.333333333333
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -83.333333333333
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 300
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 36.363636363636
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 33
$ws.Range("K15").Value = 3.030303030303
$ws.Range("L15").Value = 3.030303030303
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 3.030303030303
$ws.Range("C16").Value = 26
$ws.Range("D16").Value = 35
$ws.Range("E16").Value = -25.714285714285
$ws.Range("F16").Value = 144
$ws.Range("G16").Value = 143
$ws.Range("H16").Value = 0.6993006993
$ws.Range("I16").Value = 328
$ws.Range("J16").Value = 268
$ws.Range("K16").Value = 22.388059701492
$ws.Range("L16").Value = 87.428571428571
$ws.Range("M16").Value = -5.475504322766
$ws.Range("N16").Value = -81.007527504342
$ws.Range("C17").Value = 51
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 41.666666666666
$ws.Range("F17").Value = 189
$ws.Range("G17").Value = 151
$ws.Range("H17").Value = 25.165562913907
$ws.Range("I17").Value = 438
$ws.Range("J17").Value = 372
$ws.Range("K17").Value = 17.741935483871
$ws.Range("L17").Value = 58.12274368231
$ws.Range("M17").Value = 102.777777777778
$ws.Range("N17").Value = 0.921658986175
$ws.Range("C18").Value = 40
$ws.Range("D18").Value = 35
$ws.Range("E18").Value = 14.285714285714
$ws.Range("F18").Value = 175
$ws.Range("G18").Value = 156
$ws.Range("H18").Value = 12.179487179487
$ws.Range("I18").Value = 391
$ws.Range("J18").Value = 345
$ws.Range("K18").Value = 13.333333333333
$ws.Range("L18").Value = 22.1875
$ws.Range("M18").Value = -25.665399239543
$ws.Range("N18").Value = -86.173974540311
$ws.Range("C19").Value = 146
$ws.Range("D19").Value = 131
$ws.Range("E19").Value = 11.450381679389
$ws.Range("F19").Value = 526
$ws.Range("G19").Value = 543
$ws.Range("H19").Value = -3.130755064456
$ws.Range("I19").Value = 1152
$ws.Range("J19").Value = 1397
$ws.Range("K19").Value = -17.537580529706
$ws.Range("L19").Value = 84.32
$ws.Range("M19").Value = 77.503852080123
$ws.Range("N19").Value = -13.448534936138
$ws.Range("C20").Value = 41
$ws.Range("D20").Value = 22
$ws.Range("E20").Value = 86.363636363636
$ws.Range("F20").Value = 172
$ws.Range("G20").Value = 100
$ws.Range("H20").Value = 72
$ws.Range("I20").Value = 380
$ws.Range("J20").Value = 244
$ws.Range("K20").Value = 55.737704918032
$ws.Range("L20").Value = 126.190476190476
$ws.Range("M20").Value = 36.690647482014
$ws.Range("N20").Value = -91.146318732525
$ws.Range("C21").Value = 308
$ws.Range("D21").Value = 260
$ws.Range("E21").Value = 18.461538461538
$ws.Range("F21").Value = 1223
$ws.Range("G21").Value = 1107
$ws.Range("H21").Value = 10.478771454381
$ws.Range("I21").Value = 2727
$ws.Range("J21").Value = 2665
$ws.Range("K21").Value = 2.326454033771
$ws.Range("L21").Value = 70.4375
$ws.Range("M21").Value = 33.676470588235
$ws.Range("N21").Value = -74.439966257381
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = 16.666666666666
$ws.Range("I22").Value = 50
$ws.Range("J22").Value = 51
$ws.Range("K22").Value = -1.960784313725
$ws.Range("L22").Value = 455.555555555556
$ws.Range("M22").Value = 38.888888888888
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -21.052631578947
$ws.Range("I23").Value = 39
$ws.Range("J23").Value = 42
$ws.Range("K23").Value = -7.142857142857
$ws.Range("L23").Value = 25.806451612903
$ws.Range("M23").Value = 34.482758620689
$ws.Range("C24").Value = 307
$ws.Range("D24").Value = 292
$ws.Range("E24").Value = 5.136986301369
$ws.Range("F24").Value = 1217
$ws.Range("G24").Value = 1148
$ws.Range("H24").Value = 6.010452961672
$ws.Range("I24").Value = 2706
$ws.Range("J24").Value = 2509
$ws.Range("K24").Value = 7.851733758469
$ws.Range("L24").Value = 38.983050847457
$ws.Range("M24").Value = 89.894736842105
$ws.Range("C25").Value = 79
$ws.Range("D25").Value = 94
$ws.Range("E25").Value = -15.95744680851
$ws.Range("F25").Value = 344
$ws.Range("G25").Value = 386
$ws.Range("H25").Value = -10.880829015544
$ws.Range("I25").Value = 824
$ws.Range("J25").Value = 811
$ws.Range("K25").Value = 1.602959309494
$ws.Range("L25").Value = 38.255033557047
$ws.Range("M25").Value = 3.909205548549
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 51
$ws.Range("J26").Value = 45
$ws.Range("K26").Value = 13.333333333333
$ws.Range("L26").Value = 13.333333333333
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = 49
$ws.Range("G27").Value = 43
$ws.Range("H27").Value = 13.953488372093
$ws.Range("I27").Value = 104
$ws.Range("J27").Value = 87
$ws.Range("K27").Value = 19.540229885057
$ws.Range("L27").Value = 46.478873239436
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 37.5
$ws.Range("L28").Value = 10
$ws.Range("M28").Value = 57.142857142857
$ws.Range("N28").Value = -77.551020408163
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 25
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = 25
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 66.666666666666
$ws.Range("N29").Value = -79.591836734693
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = -71.428571428571
$ws.Range("I30").Value = 8
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = -52.941176470588
$ws.Range("L30").Value = 60
